$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 11 ("Jurisdiction" property), pushing every
# row below it down by one. Copy formatting from the row that will become
# row 12 ("Description") so the inserted row keeps the same style (s="2").
$ws.Rows(11).Insert()
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Populate the newly inserted row.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Update the Date value (row 8).
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Update the Contact value (row 10).
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# The Description row, now at row 12, keeps its original value.
$ws.Range("B12").Value = "Full Size Extension"

# Append a new row 21 ("Context" property), matching the new last metadata
# entry. Copy formatting from row 20 first.
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Context"
$ws.Range("B21").Value = "element:Element"
